# Add a new "Estudio" column (X) to the samples listing, so samples can be
# associated with a study from the Excel sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from the last existing header cell (W1) onto the new
# header cell (X1), then set its text.
$ws.Range("W1").Copy() | Out-Null
$ws.Range("X1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("X1").Value = "Estudio"

$excel.CutCopyMode = $false

# Match the author's final selection state.
$ws.Range("X1").Select() | Out-Null
